# Junction_Flooding_208 - custom accuracy + 1000 data points
# 1) Row 5 numeric values are re-rounded to the new "custom accuracy" (2 dp)
# 2) Row 6 (the extra reading) is removed entirely
# 3) A handful of data columns shrink from width 8 to width 7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update row 5 values to the reduced-precision readings ---
$row5 = @{
    "B5"  = 20.55
    "C5"  = 15.28
    "D5"  = 1.24
    "E5"  = 44.93
    "F5"  = 36.53
    "G5"  = 16.12
    "H5"  = 60.79
    "I5"  = 25.01
    "J5"  = 11.07
    "K5"  = 16.33
    "L5"  = 18.01
    "M5"  = 19.17
    "N5"  = 5.19
    "O5"  = 16.17
    "P5"  = 22.97
    "Q5"  = 13.68
    "R5"  = 0.79
    "S5"  = 0.84
    "T5"  = 238.75
    "U5"  = 45.17
    "V5"  = 14.92
    "W5"  = 30.3
    "X5"  = 15.89
    "Y5"  = 2.42
    "Z5"  = 30.03
    "AA5" = 13.18
    "AB5" = 11.71
    "AC5" = 13.77
    "AD5" = 18.91
    "AE5" = 0.54
    "AF5" = 55.34
    "AG5" = 8.369999999999999
    "AH5" = 18.66
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# --- 2) Drop the last data row (row 6) completely ---
$ws.Rows.Item(6).Delete()

# --- 3) Narrow the listed columns from 8 to 7 character units ---
# (ColumnWidth is expressed in Excel's character-width units, which are
#  offset from the raw OOXML <col width> by ~5/6; 6.1666... round-trips to 7)
$narrowColumns = @("C", "G", "J", "K", "Q", "V", "X", "AA", "AB", "AC")
foreach ($col in $narrowColumns) {
    $ws.Range("$($col):$($col)").ColumnWidth = 6.166666666666667
}
